# Testing.docx edit script
# Implements:
#  1. Paragraph "3.0" text change (add " marker" before comma, "left" -> "right")
#  2. Insert two brand new paragraphs "5.0"/"6.0" about marker colours after para "4.0"
#  3. Renumber the remaining numbered paragraphs (old 4..14 -> new 7..17)
#  4. Drop the trailing blank heading paragraph + red "_GoBack" paragraph,
#     keeping the final tab-only paragraph and re-homing the bookmark onto it.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Paragraph 5 ("3.0 To place a marker ...") - reword
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("an accident, which are located on the top left hand", $true, $false, $false, $false, $false, $true, 1, $false, "an accident marker, which are located on the top right hand", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert the two new paragraphs after paragraph 6 ("4.0 When you select...")
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs.Item(6)
$anchor.Range.InsertParagraphAfter() | Out-Null

$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "5.0" + [string][char]9 + "Accident markers are represented by a yellow " + [string][char]8220 + "A" + [string][char]8221 + " marker."

$p7b = $d.Paragraphs.Item(7)
$p7b.Range.InsertParagraphAfter() | Out-Null

$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "6.0" + [string][char]9 + "Potholes are represented by an orange " + [string][char]8220 + "P" + [string][char]8221 + " marker."

# ---------------------------------------------------------------------------
# 3. Renumber paragraphs 9..19 (old numbers 4..14 -> new numbers 7..17)
#    Each paragraph's ordinal is the first run of its Range; replace just
#    that leading substring so the tab / body runs are left untouched.
# ---------------------------------------------------------------------------
$renumbers = @(
    @{ Index = 9;  Old = "4";  New = "7" },
    @{ Index = 10; Old = "5";  New = "8" },
    @{ Index = 11; Old = "6";  New = "9" },
    @{ Index = 12; Old = "7";  New = "10" },
    @{ Index = 13; Old = "8";  New = "11" },
    @{ Index = 14; Old = "9";  New = "12" },
    @{ Index = 15; Old = "10"; New = "13" },
    @{ Index = 16; Old = "11"; New = "14" },
    @{ Index = 17; Old = "12"; New = "15" },
    @{ Index = 18; Old = "13"; New = "16" },
    @{ Index = 19; Old = "14"; New = "17" }
)

foreach ($item in $renumbers) {
    $para = $d.Paragraphs.Item($item.Index)
    $start = $para.Range.Start
    $len = $item.Old.Length
    $numRange = $d.Range($start, $start + $len)
    $numRange.Text = $item.New
}

# ---------------------------------------------------------------------------
# 4. Drop the trailing heading-style blank paragraph and the red bookmark
#    paragraph, keeping the tab-only paragraph before them. Re-home the
#    _GoBack bookmark onto the surviving paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Bookmarks.Item("_GoBack").Delete()
$tabPara = $d.Paragraphs.Item($d.Paragraphs.Count - 2)
$tabEnd = $tabPara.Range.End
$bookmarkRange = $d.Range($tabEnd - 1, $tabEnd - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$headingPara.Range.Delete() | Out-Null
$redPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$redPara.Range.Delete() | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
